# Commit: "Ajout log sur carte SD" -- add a Spi/SD component box on the
# architecture diagram (slide 1) and a matching row in the component
# table (slide 2), and slide the "SerialDriver" box left to make room.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 1 : architecture diagram
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

# Locate the "SerialDriver" rounded-rectangle box.
$serialDriver = $null
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $candidate = $s1.Shapes.Item($i)
    if ($candidate.TextFrame.TextRange.Text -eq "SerialDriver") {
        $serialDriver = $candidate
    }
}

# Slide "SerialDriver" a bit to the left to free space for the new box.
$serialDriver.Left = 2115464 / 12700
$serialDriver.Top = 889173 / 12700

# Duplicate it to get an identical style (fill/line/font theme refs)
# for the new "Spi/SD" box, then move, resize and relabel the copy.
$dup = $serialDriver.Duplicate()
$spiSd = $dup.Item(1)
$spiSd.Left = 3708466 / 12700
$spiSd.Top = 878657 / 12700
$spiSd.Width = 661603 / 12700
$spiSd.Height = 643233 / 12700
$spiSd.TextFrame.TextRange.Text = "Spi/SD"

# ---------------------------------------------------------------------
# Slide 2 : component responsibility table
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tblShape = $s2.Shapes.Item(1)
$tbl = $tblShape.Table

# Find the "GpioTools" row so the new "Spi/SD" row is inserted right
# after it (and therefore right before the "Robot2017" row).
$insertAt = $tbl.Rows.Count + 1
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    if ($tbl.Rows.Item($i).Cells.Item(1).Shape.TextFrame.TextRange.Text -eq "GpioTools") {
        $insertAt = $i + 1
    }
}

$tbl.Rows.Add($insertAt)
$newRow = $tbl.Rows.Item($insertAt)

$newRow.Cells.Item(1).Shape.TextFrame.TextRange.Text = "Spi/SD"
$newRow.Cells.Item(1).Shape.TextFrame.TextRange.Font.Size = 14

$newRow.Cells.Item(2).Shape.TextFrame.TextRange.Text = "Manage Spi bus, read SD cards, manage filesystem."
$newRow.Cells.Item(2).Shape.TextFrame.TextRange.Font.Size = 14

$newRow.Height = 290333 / 12700

# The table placeholder grows to accommodate the extra row.
$tblShape.Height = 6503496 / 12700
